{"js": "const replacements = [\n  [\"2025-10-12 Sunday\", \"2025-10-13 Monday\"],\n  [\"73\u00d780=5840\", \"98\u00d734=3332\"],\n  [\"69\u00d761=4209\", \"59\u00d765=3835\"],\n  [\"26\u00d721=546\", \"28\u00d755=1540\"],\n  [\"21\u00d742=882\", \"26\u00d778=2028\"],\n  [\"51\u00d741=2091\", \"90\u00d748=4320\"],\n  [\"85\u00d766=5610\", \"21\u00d754=1134\"],\n  [\"83\u00d742=3486\", \"18\u00d743=774\"],\n  [\"31\u00d761=1891\", \"72\u00d734=2448\"],\n  [\"92\u00d750=4600\", \"31\u00d711=341\"],\n  [\"86\u00d753=4558\", \"40\u00d757=2280\"],\n  [\"26\u00d732=832\", \"48\u00d798=4704\"],\n  [\"46\u00d787=4002\", \"70\u00d792=6440\"],\n  [\"57\u00d768=3876\", \"64\u00d763=4032\"],\n  [\"63\u00d720=1260\", \"21\u00d792=1932\"],\n  [\"71\u00d775=5325\", \"45\u00d780=3600\"],\n  [\"11\u00d736=396\", \"67\u00d759=3953\"],\n  [\"27\u00d738=1026\", \"35\u00d779=2765\"],\n  [\"93\u00d770=6510\", \"71\u00d790=6390\"],\n  [\"62\u00d764=3968\", \"71\u00d760=4260\"],\n  [\"23\u00d774=1702\", \"92\u00d754=4968\"],\n  [\"27\u00d740=1080\", \"97\u00d723=2231\"],\n  [\"85\u00d724=2040\", \"79\u00d798=7742\"],\n  [\"56\u00d767=3752\", \"23\u00d785=1955\"],\n  [\"20\u00d731=620\", \"36\u00d753=1908\"],\n  [\"30\u00d767=2010\", \"79\u00d713=1027\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-12 Sunday\", \"2025-10-13 Monday\"),\n    @(\"73\u00d780=5840\", \"98\u00d734=3332\"),\n    @(\"69\u00d761=4209\", \"59\u00d765=3835\"),\n    @(\"26\u00d721=546\", \"28\u00d755=1540\"),\n    @(\"21\u00d742=882\", \"26\u00d778=2028\"),\n    @(\"51\u00d741=2091\", \"90\u00d748=4320\"),\n    @(\"85\u00d766=5610\", \"21\u00d754=1134\"),\n    @(\"83\u00d742=3486\", \"18\u00d743=774\"),\n    @(\"31\u00d761=1891\", \"72\u00d734=2448\"),\n    @(\"92\u00d750=4600\", \"31\u00d711=341\"),\n    @(\"86\u00d753=4558\", \"40\u00d757=2280\"),\n    @(\"26\u00d732=832\", \"48\u00d798=4704\"),\n    @(\"46\u00d787=4002\", \"70\u00d792=6440\"),\n    @(\"57\u00d768=3876\", \"64\u00d763=4032\"),\n    @(\"63\u00d720=1260\", \"21\u00d792=1932\"),\n    @(\"71\u00d775=5325\", \"45\u00d780=3600\"),\n    @(\"11\u00d736=396\", \"67\u00d759=3953\"),\n    @(\"27\u00d738=1026\", \"35\u00d779=2765\"),\n    @(\"93\u00d770=6510\", \"71\u00d790=6390\"),\n    @(\"62\u00d764=3968\", \"71\u00d760=4260\"),\n    @(\"23\u00d774=1702\", \"92\u00d754=4968\"),\n    @(\"27\u00d740=1080\", \"97\u00d723=2231\"),\n    @(\"85\u00d724=2040\", \"79\u00d798=7742\"),\n    @(\"56\u00d767=3752\", \"23\u00d785=1955\"),\n    @(\"20\u00d731=620\", \"36\u00d753=1908\"),\n    @(\"30\u00d767=2010\", \"79\u00d713=1027\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
